# Remove the "esim_provider" column (column L) entirely.
# This shifts "remarks" (previously column M) left into column L,
# matching the target workbook layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(12).Delete()

$ws.Range("S10").Select()
